$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7139005064964294
$ws.Range("B1").Value = 0.7739047408103943
$ws.Range("C1").Value = 0.8943632245063782
$ws.Range("D1").Value = 1.460410118103027
$ws.Range("E1").Value = 4.20029878616333
